$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# Resize column H (8) so its stored width becomes 17.5
$ws.Columns.Item(8).ColumnWidth = 16.6625

# Add "Priority" (column G) values for rows 2-12
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 1

# Add "Effort Actual [h]" (column I) values for rows 2-3
$ws.Range("I2").Value = 2
$ws.Range("I3").Value = 4

# Update the selection to D3 (active cell)
$ws.Range("D3").Select()
